$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "58.536.95"
$ws.Range("E2").Value = "  -4.91%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.462.50"
$ws.Range("E3").Value = "  -4.51%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "532.53"
$ws.Range("E5").Value = "  -3.67%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.58"
$ws.Range("E6").Value = "  -7.12%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  -0.25%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.566"
$ws.Range("E8").Value = "  -5.15%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.481.08"
$ws.Range("E9").Value = "  -3.91%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0991"
$ws.Range("E10").Value = "  -4.78%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.53"
$ws.Range("E12").Value = "  +1.50%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.350"
$ws.Range("E13").Value = "  -3.58%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.900.71"
$ws.Range("E14").Value = "  -4.46%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "23.65"
$ws.Range("E15").Value = "  -6.87%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "58.473.57"
$ws.Range("E16").Value = "  -4.89%  "
$ws.Range("E17").Value = "  -4.50%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.482.67"
$ws.Range("E18").Value = "  -3.87%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.23"
$ws.Range("E19").Value = "  -2.81%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.29"
$ws.Range("E20").Value = "  -5.16%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "321.77"
$ws.Range("E21").Value = "  -4.83%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.997"
$ws.Range("E22").Value = "  -0.19%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.71"
$ws.Range("E23").Value = "  -5.41%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "60.57"
$ws.Range("E24").Value = "  -3.64%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.436"
$ws.Range("E25").Value = "  -11.57%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.996"
$ws.Range("E26").Value = "  -0.22%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.160"
$ws.Range("E27").Value = "  -4.77%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.579.27"
$ws.Range("E28").Value = "  -4.59%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.66"
$ws.Range("E29").Value = "  -4.81%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.89"
$ws.Range("E30").Value = "  -2.08%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0₃0770"
$ws.Range("E31").Value = "  -7.97%  "
$ws.Range("E32").Value = "  -6.65%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.22"
$ws.Range("E33").Value = "  -5.85%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.998"
$ws.Range("E34").Value = "  -0.08%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "157.61"
$ws.Range("E35").Value = "  -1.65%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.39"
$ws.Range("E36").Value = "  -1.29%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "18.41"
$ws.Range("E37").Value = "  -4.01%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.37"
$ws.Range("E38").Value = "  -6.24%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.59"
$ws.Range("E39").Value = "  -10.67%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.71"
$ws.Range("E40").Value = "  -5.05%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "303.78"
$ws.Range("E41").Value = "  -9.67%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "36.44"
$ws.Range("E42").Value = "  -2.54%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.69"
$ws.Range("E43").Value = "  -6.10%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.800"
$ws.Range("E44").Value = "  -10.27%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.995"
$ws.Range("E45").Value = "  -0.19%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.77"
$ws.Range("E46").Value = "  -1.41%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.588"
$ws.Range("E47").Value = "  -2.80%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "123.83"
$ws.Range("E48").Value = "  -0.17%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0921"
$ws.Range("E49").Value = "  -4.57%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0517"
$ws.Range("E50").Value = "  -4.97%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0227"
$ws.Range("E51").Value = "  -5.10%  "
